# "submission groups.xlsx" has 7 sheets named "11".."17", one per
# submission round. The roster of names for the current round lives on
# sheet "17" (its A1:A4), and that sheet is the active/selected tab.
#
# This upload starts a new round: the four names are replaced with a new
# group of students, and the roster (now for the new round) is moved onto
# sheet "11", which becomes the active/selected tab again; sheet "17" goes
# back to being empty, just like sheet "11" used to be.

$wb = $excel.ActiveWorkbook

$sheet11 = $wb.Worksheets.Item("11")
$sheet17 = $wb.Worksheets.Item("17")

# Rename the four names in place on sheet "17" first (this reuses the
# same shared-string slots that already back these four cells).
$sheet17.Range("A1").Value = "batel elbaz"
$sheet17.Range("A2").Value = "shahar gavriel"
$sheet17.Range("A3").Value = "liad tzvaot"
$sheet17.Range("A4").Value = "idan yontov"

# Write the same four names onto sheet "11", in the new row order.
$sheet11.Range("A1").Value = "batel elbaz"
$sheet11.Range("A2").Value = "shahar gavriel"
$sheet11.Range("A3").Value = "idan yontov"
$sheet11.Range("A4").Value = "liad tzvaot"

# The roster no longer lives on sheet "17" - clear it back out.
$sheet17.Range("A1:A4").ClearContents()

# Sheet "11" becomes the active / selected tab (it was "17" before).
[void]$sheet11.Activate()
[void]$sheet11.Range("E3").Select()

# Sheet "17" just keeps a plain cell selection now that it's empty again.
[void]$sheet17.Range("D7").Select()

# Leave sheet "11" as the active sheet.
[void]$sheet11.Activate()
